$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.759736666666667
$ws.Range("H2").Value = 11.27921
$ws.Range("I2").Value = 0.0683751702595819
$ws.Range("J2").Value = 0.06837517025958188
$ws.Range("M2").Value = 1.363346333333333
$ws.Range("N2").Value = 4.090039
$ws.Range("O2").Value = 0.02430403345239443
$ws.Range("P2").Value = 0.02430403345239443
$ws.Range("Q2").Value = 5.125823198798889
$ws.Range("R2").Value = 46.13240878918999
$ws.Range("S2").Value = 0.001661792425302043
$ws.Range("T2").Value = 0.001661792425302043

# Row 3
$ws.Range("G3").Value = 3.759736666666667
$ws.Range("H3").Value = 11.27921
$ws.Range("I3").Value = 0.0683751702595819
$ws.Range("J3").Value = 0.06837517025958188
$ws.Range("O3").Value = 0.679596855668023
$ws.Range("P3").Value = 0.679596855668023
$ws.Range("Q3").Value = 143.3298442185422
$ws.Range("R3").Value = 1289.96859796688
$ws.Range("S3").Value = 0.04646755071417757
$ws.Range("T3").Value = 0.04646755071417757

# Row 4
$ws.Range("G4").Value = 3.759736666666667
$ws.Range("H4").Value = 11.27921
$ws.Range("I4").Value = 0.0683751702595819
$ws.Range("J4").Value = 0.06837517025958188
$ws.Range("N4").Value = 49.82946200000001
$ws.Range("O4").Value = 0.2960991108795826
$ws.Range("P4").Value = 0.2960991108795826
$ws.Range("Q4").Value = 62.44855178722445
$ws.Range("R4").Value = 562.03696608502
$ws.Range("S4").Value = 0.02024582712010228
$ws.Range("T4").Value = 0.02024582712010228

# Row 5
$ws.Range("I5").Value = 0.6514180024294648
$ws.Range("J5").Value = 0.6514180024294647
$ws.Range("M5").Value = 1.363346333333333
$ws.Range("N5").Value = 4.090039
$ws.Range("O5").Value = 0.02430403345239443
$ws.Range("P5").Value = 0.02430403345239443
$ws.Range("Q5").Value = 48.83429900491188
$ws.Range("R5").Value = 439.508691044207
$ws.Range("S5").Value = 0.01583208492253767
$ws.Range("T5").Value = 0.01583208492253767

# Row 6
$ws.Range("I6").Value = 0.6514180024294648
$ws.Range("J6").Value = 0.6514180024294647
$ws.Range("O6").Value = 0.679596855668023
$ws.Range("P6").Value = 0.679596855668023
$ws.Range("S6").Value = 0.4427016261766088
$ws.Range("T6").Value = 0.4427016261766087

# Row 7
$ws.Range("I7").Value = 0.6514180024294648
$ws.Range("J7").Value = 0.6514180024294647
$ws.Range("N7").Value = 49.82946200000001
$ws.Range("O7").Value = 0.2960991108795826
$ws.Range("P7").Value = 0.2960991108795826
$ws.Range("S7").Value = 0.1928842913303183
$ws.Range("T7").Value = 0.1928842913303183

# Row 8
$ws.Range("I8").Value = 0.2802068273109533
$ws.Range("J8").Value = 0.2802068273109533
$ws.Range("M8").Value = 1.363346333333333
$ws.Range("N8").Value = 4.090039
$ws.Range("O8").Value = 0.02430403345239443
$ws.Range("P8").Value = 0.02430403345239443
$ws.Range("Q8").Value = 21.00602675561222
$ws.Range("R8").Value = 189.05424080051
$ws.Range("S8").Value = 0.006810156104554718
$ws.Range("T8").Value = 0.006810156104554717

# Row 9
$ws.Range("I9").Value = 0.2802068273109533
$ws.Range("J9").Value = 0.2802068273109533
$ws.Range("O9").Value = 0.679596855668023
$ws.Range("P9").Value = 0.679596855668023
$ws.Range("S9").Value = 0.1904276787772366
$ws.Range("T9").Value = 0.1904276787772365

# Row 10
$ws.Range("I10").Value = 0.2802068273109533
$ws.Range("J10").Value = 0.2802068273109533
$ws.Range("N10").Value = 49.82946200000001
$ws.Range("O10").Value = 0.2960991108795826
$ws.Range("P10").Value = 0.2960991108795826
$ws.Range("S10").Value = 0.08296899242916202
$ws.Range("T10").Value = 0.082968992429162
